$wb = $excel.ActiveWorkbook
$wsDash = $wb.Worksheets.Item("dashboard")
$ws = $wb.Worksheets.Item("writing")

# --- Append the new day's row (row 26) to the "writing" log ---------------
$ws.Range("A26").Value = 44159
$ws.Range("B26").Value = 236
$ws.Range("C26").Value = 87
$ws.Range("D26").Value = 516
$ws.Range("E26").Value = 7486
$ws.Range("F26").Value = 7980
$ws.Range("G26").Value = 120
$ws.Range("H26").Value = 117
$ws.Range("I26").Value = 186
$ws.Range("J26").Formula = "=SUM(B26:I26)"
$ws.Range("K26").Formula = "=J26-J25"

# Match the date formatting used by the rest of column A (copy format only)
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Grow Table1 so the new row is part of the structured table -----------
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:K26"))

# --- Extend the dashboard chart's series ranges to include the new row ----
$co = $wsDash.ChartObjects().Item(1)
$chart = $co.Chart
$sDaily = $chart.SeriesCollection().Item(1)
$sTotal = $chart.SeriesCollection().Item(2)
$sDaily.Formula = "=SERIES(writing!`$K`$1,writing!`$A`$2:`$A`$26,writing!`$K`$2:`$K`$26,1)"
$sTotal.Formula = "=SERIES(writing!`$J`$1,writing!`$A`$2:`$A`$26,writing!`$J`$2:`$J`$26,2)"

# --- Leave the selection on the new row, but keep "dashboard" the active tab
$ws.Range("E26").Select()
$wsDash.Activate()
